# corrected ICDC Breed 1-14 scripts
#
# The FilesTab Cypher query (cell B4 on the "startup" sheet) dropped the
# `File Type` and `Breed` columns from its RETURN clause. Updating the
# cell text also reflows the wrapped row height, and the active selection
# on the sheet moved from A4 to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Dalmatian']  
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# Update the query text in B4 (removes the File Type / Breed columns).
$ws.Range("B4").Value = $newQuery

# The shorter text re-wraps to a smaller row height (matches row 3's height).
$ws.Rows.Item(4).RowHeight = 217.5

# Selection moved from A4 to B4.
$ws.Activate() | Out-Null
$ws.Range("B4").Select() | Out-Null
